$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.696.37'
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("D3").Value = '3.571.19'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.93'
$ws.Range("E5").Value = '  +4.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.89'
$ws.Range("E6").Value = '  +0.47%  '
$ws.Range("D7").Value = '3.565.05'
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.197'
$ws.Range("E10").Value = '  +3.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.42'
$ws.Range("E11").Value = '  +9.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.589'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.89'
$ws.Range("E13").Value = '  -1.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000278'
$ws.Range("E14").Value = '  +0.60%  '
$ws.Range("D15").Value = '4.148.72'
$ws.Range("E15").Value = '  +0.92%  '
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '616.49'
$ws.Range("E17").Value = '  -2.27%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '70.793.65'
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.553.00'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("E20").Value = '  -2.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.44'
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.889'
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.37'
$ws.Range("E23").Value = '  -16.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.04'
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.53'
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.57'
$ws.Range("E29").Value = '  +1.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.13'
$ws.Range("E30").Value = '  -2.58%  '
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("E32").Value = '  -4.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.03'
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("E34").Value = '  -2.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '613.18'
$ws.Range("E35").Value = '  -4.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.78'
$ws.Range("E36").Value = '  +6.98%  '
$ws.Range("E37").Value = '  -1.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.86'
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0482'
$ws.Range("E39").Value = '  +5.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '57.34'
$ws.Range("E40").Value = '  -0.28%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  +3.55%  '
$ws.Range("D43").Value = '3.386.78'
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("E44").Value = '  -3.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.01'
$ws.Range("E45").Value = '  +8.43%  '
$ws.Range("D46").Value = '0.0₃0711'
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '33.02'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("E50").Value = '  -0.20%  '

Write-Host "Applied crypto list update"
